$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows newly marked "gescaled" (scaled) = yes, and assigned to "lara" (new hire),
# with level (E) = 1 where indicated by the diff.
$ws.Range("C13").Value = "yes"
$ws.Range("D13").Value = "lara"

$ws.Range("C22").Value = "yes"
$ws.Range("D22").Value = "lara"
$ws.Range("E22").Value = 1

$ws.Range("C23").Value = "yes"
$ws.Range("D23").Value = "lara"
$ws.Range("E23").Value = 1

# Row 30: "staat het op git?" flips from Bad(red) to Good(green); gescaled = yes
$ws.Range("B30").Style = "Good"
$ws.Range("C30").Value = "yes"

# Row 43: gescaled = yes, level = 1 (already assigned to marc)
$ws.Range("C43").Value = "yes"
$ws.Range("E43").Value = 1

$ws.Range("C47").Value = "yes"
$ws.Range("D47").Value = "lara"
$ws.Range("E47").Value = 1

$ws.Range("C48").Value = "yes"
$ws.Range("D48").Value = "lara"
$ws.Range("E48").Value = 1

$ws.Range("C49").Value = "yes"
$ws.Range("D49").Value = "lara"
$ws.Range("E49").Value = 1

$ws.Range("C50").Value = "yes"
$ws.Range("D50").Value = "lara"
$ws.Range("E50").Value = 1

$ws.Range("C51").Value = "yes"
$ws.Range("D51").Value = "lara"
$ws.Range("E51").Value = 1

# Row 54: gescaled = yes (already assigned to marc)
$ws.Range("C54").Value = "yes"

$ws.Range("C58").Value = "yes"
$ws.Range("D58").Value = "lara"

# Row 59: "staat het op git?" flips from Bad to Good; gescaled = yes (already assigned to marc)
$ws.Range("B59").Style = "Good"
$ws.Range("C59").Value = "yes"

# Row 65: "staat het op git?" flips from Bad to Good; gescaled = yes; assigned to lara
$ws.Range("B65").Style = "Good"
$ws.Range("C65").Value = "yes"
$ws.Range("D65").Value = "lara"

# Rows 67, 72, 73, 76, 78: gescaled = yes (already assigned to marc)
$ws.Range("C67").Value = "yes"
$ws.Range("C72").Value = "yes"
$ws.Range("C73").Value = "yes"
$ws.Range("C76").Value = "yes"
$ws.Range("C78").Value = "yes"

# Row 79: "staat het op git?" becomes Bad(red, newly tracked); gescaled = yes (already assigned to marc)
$ws.Range("B79").Style = "Bad"
$ws.Range("C79").Value = "yes"

# Move the view / current selection to reflect where work left off.
$ws.Range("C78").Select()
